# Folder-structure refactor: the "rbwatson" path segment used throughout the
# WLUX / Localhost endpoint formulas (and the embedded JSON sample in the
# shared strings table) is renamed to "data".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Activate()

# Replace every occurrence of "rbwatson" with "data" across the formulas /
# text on Sheet1 (this covers the CONCATENATE(...) endpoint-builder formulas
# in columns F and G for rows 2-17).
$ws.Cells.Replace("rbwatson", "data")

# Leave the selection where the author left it after the refactor.
$ws.Range("G2").Select()
